# NYPD 112th Precinct CompStat weekly report - roll forward to the new reporting week.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number and the "week covering" date range ---
$ws.Range("A8").Value = "Volume 32   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/10/2025  Through  11/16/2025"

# --- Row 14 ---
$ws.Range("N14").Value = -60

# --- Row 15 ---
$ws.Range("M15").Value = 75

# --- Row 16 ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 67
$ws.Range("J16").Value = 77
$ws.Range("K16").Value = -12.987012987013
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -37.383177570093
$ws.Range("N16").Value = -88.224956063268

# --- Row 17 ---
$ws.Range("C17").Value = 5
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0"
$ws.Range("E17").Value = "***.*"
$ws.Range("F17").Value = 13
$ws.Range("H17").Value = 62.5
$ws.Range("I17").Value = 139
$ws.Range("K17").Value = 41.836734693877
$ws.Range("L17").Value = 54.444444444444
$ws.Range("M17").Value = 172.549019607843
$ws.Range("N17").Value = 23.008849557522

# --- Row 18 ---
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 50
$ws.Range("I18").Value = 145
$ws.Range("J18").Value = 89
$ws.Range("K18").Value = 62.921348314606
$ws.Range("L18").Value = 42.156862745098
$ws.Range("M18").Value = 29.464285714285
$ws.Range("N18").Value = -88.409272581934

# --- Row 19 ---
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 133.333333333333
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 29.032258064516
$ws.Range("I19").Value = 355
$ws.Range("J19").Value = 403
$ws.Range("K19").Value = -11.910669975186
$ws.Range("L19").Value = -15.274463007159
$ws.Range("M19").Value = 4.105571847507
$ws.Range("N19").Value = -59.932279909706

# --- Row 20 ---
$ws.Range("C20").Value = 6
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 7.142857142857
$ws.Range("I20").Value = 149
$ws.Range("J20").Value = 152
$ws.Range("K20").Value = -1.973684210526
$ws.Range("L20").Value = -1.973684210526
$ws.Range("M20").Value = 69.318181818181
$ws.Range("N20").Value = -95.053120849933

# --- Row 21 ---
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 123.076923076923
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = 52.380952380952
$ws.Range("I21").Value = 864
$ws.Range("J21").Value = 830
$ws.Range("K21").Value = 4.096385542168
$ws.Range("L21").Value = 3.473053892215
$ws.Range("M21").Value = 22.727272727272
$ws.Range("N21").Value = -85.218135158254

# --- Row 22 ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 25
$ws.Range("K22").Value = -24.242424242424
$ws.Range("L22").Value = -3.846153846153
$ws.Range("M22").Value = 19.047619047619

# --- Row 24 ---
$ws.Range("C24").Value = 54
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 205
$ws.Range("G24").Value = 147
$ws.Range("H24").Value = 39.455782312925
$ws.Range("I24").Value = 1760
$ws.Range("J24").Value = 1536
$ws.Range("K24").Value = 14.583333333333
$ws.Range("L24").Value = 29.126925898752
$ws.Range("M24").Value = 104.176334106729

# --- Row 25 ---
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = 33
$ws.Range("E25").Value = 12.121212121212
$ws.Range("F25").Value = 136
$ws.Range("H25").Value = 16.239316239316
$ws.Range("I25").Value = 1326
$ws.Range("J25").Value = 1129
$ws.Range("K25").Value = 17.449069973427
$ws.Range("L25").Value = 36.139630390143

# --- Row 26 ---
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 83.333333333333
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 45.833333333333
$ws.Range("I26").Value = 262
$ws.Range("J26").Value = 243
$ws.Range("K26").Value = 7.818930041152
$ws.Range("L26").Value = 20.737327188940
$ws.Range("M26").Value = 36.458333333333

# --- Row 28 ---
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 17.142857142857
$ws.Range("L28").Value = 17.142857142857

# --- Row 31 ---
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("H31").Value = "***.*"
